$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 70, shifting existing rows 70-114 down to 71-115.
$ws.Rows("70").Insert()

# Populate the newly inserted row 70 with the new record.
$ws.Range("A70").Value = 2
$ws.Range("B70").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C70").Value = 'Coquimbo'
$ws.Range("D70").Value = 45001
$ws.Range("E70").Value = 4
$ws.Range("F70").Value = 100112030
$ws.Range("G70").Value = 'Poroto granado'
$ws.Range("H70").Value = 'Sin especificar'
$ws.Range("I70").Value = 'Primera'
$ws.Range("J70").Value = 400
$ws.Range("K70").Value = 21000
$ws.Range("L70").Value = 23000
$ws.Range("M70").Value = 22000
$ws.Range("N70").Value = '$/malla 25 kilos'
$ws.Range("O70").Value = 'Provincia de Limarí'
$ws.Range("P70").Value = 880
$ws.Range("Q70").Value = 25
$ws.Range("R70").Value = 'Hortaliza'
